$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 35: update "Kello" (time) and "Oppimisen sisältö" (content) cells.
# B35 gets the new (longer) time-range string and picks up wrap text,
# C35 gets the new description text, G35's hours bump from 2 to 3.5.
$ws.Range("B35").Value = "9.30-11.30, 12.00-13.30, 18.15-"
$ws.Range("C35").Value = "Laatikon ja tason, ja laatikon ja pisteen törmäystarkastelu,  Kahden laatikon törmäystarkastelu, kertailua ja uudelleenlukemista"
$ws.Range("B35").WrapText = $true
$ws.Range("G35").Value = 3.5

# Row height shrinks now that B35 wraps onto fewer lines than before.
$ws.Rows(35).RowHeight = 72.5

# Move the view: scrolled one row further and selection now on D35.
[void]$ws.Range("D35").Select()
